# Generate Report for Handback
# Swap a.md / b.md rows (b.md now reported first), mark both as handed back,
# record the handback target/file columns, and stamp the handback datetime.

$wb = $excel.ActiveWorkbook

$aUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/5468c7b6325e62644e2817426ef0e7e644c44359/e2e/a.md"
$bUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/5468c7b6325e62644e2817426ef0e7e644c44359/e2e/b.md"
$cfgUrl = "https://github.com/OpenLocalizationTest/oltest/blob/5468c7b6325e62644e2817426ef0e7e644c44359/.localization-config"
$xlfZhUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3ecc35f5c235adc3f012d8fa20167cb37b76432c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$xlfDeUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3e242f27e597d2e55780db625f8f930e13823f42/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$statusHandedBack = "Handed back: in sync with en-US"
$hyperlinkColor = 15570276  # matches the workbook's existing HyperLink cell style (FF6495ED)

function Add-Link($ws, $cellRef, $url, $text) {
    $ws.Hyperlinks.Add($ws.Range($cellRef), $url, "", "", $text)
    $ws.Range($cellRef).Font.Underline = $true
    $ws.Range($cellRef).Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------------
# Overview sheet: rows 2/3 swap which file (a.md / b.md) they describe, and
# both language columns move to "handed back".
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A1").Hyperlinks.Delete()

$ws.Range("B2").Value = $statusHandedBack
$ws.Range("C2").Value = $statusHandedBack
$ws.Range("B3").Value = $statusHandedBack
$ws.Range("C3").Value = $statusHandedBack

Add-Link $ws "A2" $bUrl "b.md"
Add-Link $ws "A3" $aUrl "a.md"
Add-Link $ws "A4" $cfgUrl ".localization-config"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A1").Hyperlinks.Delete()

$ws.Range("B2").Value = $statusHandedBack
$ws.Range("B3").Value = $statusHandedBack
$ws.Range("G2").Value = "2016-03-03 02:57:47"
$ws.Range("G3").Value = "2016-03-03 02:57:47"

Add-Link $ws "A2" $bUrl "b.md"
Add-Link $ws "C2" $xlfZhUrl "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
Add-Link $ws "E2" $bUrl "b.md"
Add-Link $ws "F2" $xlfZhUrl "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

Add-Link $ws "A3" $aUrl "a.md"
Add-Link $ws "C3" $xlfZhUrl "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
Add-Link $ws "E3" $bUrl "b.md"
Add-Link $ws "F3" $xlfZhUrl "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

Add-Link $ws "A4" $cfgUrl ".localization-config"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A1").Hyperlinks.Delete()

$ws.Range("B2").Value = $statusHandedBack
$ws.Range("B3").Value = $statusHandedBack
$ws.Range("G2").Value = "2016-03-03 02:58:06"
$ws.Range("G3").Value = "2016-03-03 02:58:06"

Add-Link $ws "A2" $bUrl "b.md"
Add-Link $ws "C2" $xlfDeUrl "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
Add-Link $ws "E2" $bUrl "b.md"
Add-Link $ws "F2" $xlfDeUrl "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

Add-Link $ws "A3" $aUrl "a.md"
Add-Link $ws "C3" $xlfDeUrl "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
Add-Link $ws "E3" $bUrl "b.md"
Add-Link $ws "F3" $xlfDeUrl "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

Add-Link $ws "A4" $cfgUrl ".localization-config"

Write-Host "Report generated for handback"
